# feat: add 2022-Q1 data
#
# 1) Prepend a "2022-Q1" row to the "总计" (totals) summary sheet, shifting
#    the existing rows down and renumbering the helper index column.
# 2) Insert a new worksheet "2022-Q1" (fund-holdings detail, same layout as
#    the other quarterly sheets) right after "2021-Q4" and before "总计".
#
# NOTE: sheet object references captured *before* a new worksheet is
# inserted can end up repointed at the newly-inserted sheet once
# Worksheets.Add() runs, so the "总计" edits are done first (while its
# handle is still guaranteed valid) and the new sheet's handle is used
# immediately after creation, without being reused across the Add() call.

$wb = $excel.ActiveWorkbook

# --- 1) Prepend "2022-Q1" to the "总计" summary sheet -----------------------

$totalSheet = $wb.Worksheets.Item("总计")

$totalSheet.Rows("2:2").Insert()
$totalSheet.Range("B2:D2").ClearFormats()
$totalSheet.Range("A3").Copy($totalSheet.Range("A2"))

$totalSheet.Cells.Item(2,1).Value = 0
$totalSheet.Cells.Item(2,2).Value = "2022-Q1"
$totalSheet.Cells.Item(2,3).Value = 5
$totalSheet.Cells.Item(2,4).Value = 1.3

# Renumber the pre-existing rows' helper index column (they each shifted
# down by one row and their 0-based sequence number increases by one).
$totalSheet.Cells.Item(3,1).Value = 1
$totalSheet.Cells.Item(4,1).Value = 2
$totalSheet.Cells.Item(5,1).Value = 3
$totalSheet.Cells.Item(6,1).Value = 4
$totalSheet.Cells.Item(7,1).Value = 5

# --- 2) New "2022-Q1" sheet -------------------------------------------------

$template = $wb.Worksheets.Item("2021-Q4")
$totalSheet = $wb.Worksheets.Item("总计")

$newSheet = $wb.Worksheets.Add($totalSheet)
$newSheet.Name = "2022-Q1"

# Bring over the header row + one formatted data row from the 2021-Q4 sheet
# so fonts / borders / alignment / column-index styling match the other
# quarterly sheets, then fan the data-row style out to the remaining rows.
$template.Range("A1:H2").Copy($newSheet.Range("A1:H2"))
$template.Range("A2:H2").Copy($newSheet.Range("A3:H3"))
$template.Range("A2:H2").Copy($newSheet.Range("A4:H4"))
$template.Range("A2:H2").Copy($newSheet.Range("A5:H5"))
$template.Range("A2:H2").Copy($newSheet.Range("A6:H6"))

# Fund code / size / position columns are stored as text (e.g. "007139"
# keeps its leading zero), so force text formatting before writing values.
$newSheet.Range("B2:G6").NumberFormat = "@"

$newSheet.Cells.Item(1,2).Value = "基金代码"
$newSheet.Cells.Item(1,3).Value = "基金名称"
$newSheet.Cells.Item(1,4).Value = "基金规模"
$newSheet.Cells.Item(1,5).Value = "股票总仓位"
$newSheet.Cells.Item(1,6).Value = "仓位占比"
$newSheet.Cells.Item(1,7).Value = "持有市值(亿元)"
$newSheet.Cells.Item(1,8).Value = "仓位排名"

$fundRows = @(
    @(0, "007139", "富国民裕进取沪港深成长精选混合", "12.79", "92.21", "5.18", "0.6625", 7),
    @(1, "004424", "汇添富文体娱乐主题混合",          "18.52", "90.11", "3.14", "0.5815", 10),
    @(2, "004099", "前海开源沪港深景气行业精选灵活配置混合", "0.41", "93.07", "8.47", "0.0347", 6),
    @(3, "006205", "汇添富沪港深优势精选定期开放混合", "0.40", "93.67", "4.37", "0.0175", 8),
    @(4, "519602", "海富通大中华精选混合QDII",        "0.11", "89.68", "4.85", "0.0053", 6)
)

for ($i = 0; $i -lt $fundRows.Count; $i++) {
    $r = $i + 2
    $row = $fundRows[$i]
    $newSheet.Cells.Item($r,1).Value = $row[0]
    $newSheet.Cells.Item($r,2).Value = $row[1]
    $newSheet.Cells.Item($r,3).Value = $row[2]
    $newSheet.Cells.Item($r,4).Value = $row[3]
    $newSheet.Cells.Item($r,5).Value = $row[4]
    $newSheet.Cells.Item($r,6).Value = $row[5]
    $newSheet.Cells.Item($r,7).Value = $row[6]
    $newSheet.Cells.Item($r,8).Value = $row[7]
}
